$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 618.5
$ws.Range("I2").Value = 501
$ws.Range("J2").Value = 736
$ws.Range("K2").Value = 501
$ws.Range("L2").Value = 736
$ws.Range("M2").Value = -388
$ws.Range("N2").Value = -962

$ws.Range("H5").Value = 390.45456
$ws.Range("I5").Value = 474.14285
$ws.Range("J5").Value = 244
$ws.Range("K5").Value = 474.14285
$ws.Range("L5").Value = 244
$ws.Range("M5").Value = -359.14285
$ws.Range("N5").Value = -474

$ws.Range("H13").Value = 2003
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2003
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 2003
$ws.Range("N13").Value = -2341

$ws.Range("H75").Value = 32000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 32000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 32000
$ws.Range("N75").Value = -33872

$ws.Range("H78").Value = 32000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 32000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 96000
$ws.Range("N78").Value = -105360

$ws.Range("H132").Value = 996.6667
$ws.Range("I132").Value = 996.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2990.0001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -460.0001000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2666.6667
$ws.Range("I25").Value = 2666.6667
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2666.6667
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2264.6667

$ws.Range("H45").Value = 1585.2222
$ws.Range("I45").Value = 1585.2222
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1585.2222
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1208.2222
$ws.Range("N45").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1443
$ws.Range("I20").Value = 1384.5714
$ws.Range("J20").Value = 1647.5
$ws.Range("K20").Value = 1384.5714
$ws.Range("L20").Value = 1647.5
$ws.Range("M20").Value = -1137.5714
$ws.Range("N20").Value = -2141.5

$ws.Range("H86").Value = 2668.2222
$ws.Range("I86").Value = 1835.8334
$ws.Range("J86").Value = 4333
$ws.Range("K86").Value = 1835.8334
$ws.Range("L86").Value = 4333
$ws.Range("M86").Value = -712.8334
$ws.Range("N86").Value = -6579

$ws.Range("H89").Value = 2668.2222
$ws.Range("I89").Value = 1835.8334
$ws.Range("J89").Value = 4333
$ws.Range("K89").Value = 9179.166999999999
$ws.Range("L89").Value = 21665
$ws.Range("M89").Value = -3563.166999999999
$ws.Range("N89").Value = -32897

$ws.Range("H99").Value = 4528
$ws.Range("I99").Value = 4528
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4528
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3030
$ws.Range("N99").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 4725
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 4725
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 4725
$ws.Range("M26").Value = $null
$ws.Range("N26").Value = -5299

$ws.Range("H132").Value = 4719.737
$ws.Range("I132").Value = 4426.3887
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 13279.1661
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -10749.1661
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 13797
$ws.Range("I4").Value = 15482.321
$ws.Range("J4").Value = 1999.75
$ws.Range("K4").Value = 46446.963
$ws.Range("L4").Value = 5999.25
$ws.Range("M4").Value = -46334.963
$ws.Range("N4").Value = -6223.25

$ws.Range("H12").Value = 241.5
$ws.Range("I12").Value = 68
$ws.Range("J12").Value = 365.42856
$ws.Range("K12").Value = 204
$ws.Range("L12").Value = 1096.28568
$ws.Range("M12").Value = -31
$ws.Range("N12").Value = -1442.28568

$ws.Range("H50").Value = 456.25
$ws.Range("I50").Value = 308.83334
$ws.Range("J50").Value = 898.5
$ws.Range("K50").Value = 926.5000200000001
$ws.Range("L50").Value = 2695.5
$ws.Range("M50").Value = -445.5000200000001
$ws.Range("N50").Value = -3657.5

$ws.Range("H53").Value = 456.25
$ws.Range("I53").Value = 308.83334
$ws.Range("J53").Value = 898.5
$ws.Range("K53").Value = 926.5000200000001
$ws.Range("L53").Value = 2695.5
$ws.Range("M53").Value = -445.5000200000001
$ws.Range("N53").Value = -3657.5

$ws.Range("H109").Value = 469
$ws.Range("I109").Value = 469
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1407
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4330
$ws.Range("I70").Value = 4330
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4330
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4060
$ws.Range("N70").Value = $null

$ws.Range("H73").Value = 4330
$ws.Range("I73").Value = 4330
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4330
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3394
$ws.Range("N73").Value = $null

$ws.Range("H80").Value = 881.5
$ws.Range("I80").Value = 881.8570999999999
$ws.Range("J80").Value = 879
$ws.Range("K80").Value = 881.8570999999999
$ws.Range("L80").Value = 879
$ws.Range("M80").Value = 116.1429000000001
$ws.Range("N80").Value = -2875

$ws.Range("H83").Value = 881.5
$ws.Range("I83").Value = 881.8570999999999
$ws.Range("J83").Value = 879
$ws.Range("K83").Value = 4409.2855
$ws.Range("L83").Value = 4395
$ws.Range("M83").Value = 582.7145
$ws.Range("N83").Value = -14379

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null

$ws.Range("H126").Value = 5333
$ws.Range("I126").Value = 5999.5
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 17998.5
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -15528.5
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4671.5713
$ws.Range("I7").Value = 4646.346
$ws.Range("J7").Value = 4999.5
$ws.Range("K7").Value = 4646.346
$ws.Range("L7").Value = 4999.5
$ws.Range("M7").Value = -4534.346
$ws.Range("N7").Value = -5223.5

$ws.Range("H55").Value = 788
$ws.Range("I55").Value = 499.2857
$ws.Range("J55").Value = 1293.25
$ws.Range("K55").Value = 499.2857
$ws.Range("L55").Value = 1293.25
$ws.Range("M55").Value = -326.2857
$ws.Range("N55").Value = -1639.25

$ws.Range("H93").Value = 4000
$ws.Range("I93").Value = 3500
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 3500
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -2252
$ws.Range("N93").Value = -7496

$ws.Range("H126").Value = 4671.5713
$ws.Range("I126").Value = 4646.346
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 13939.038
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -11469.038
$ws.Range("N126").Value = -19938.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 1199.5
$ws.Range("I18").Value = 400
$ws.Range("J18").Value = 1999
$ws.Range("K18").Value = 400
$ws.Range("L18").Value = 1999
$ws.Range("M18").Value = -227
$ws.Range("N18").Value = -2345

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("N83").Value = $null
